$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 115
$ws1.Range("F4").Value = 500
$ws1.Range("F5").Value = 5063
$ws1.Range("F10").Value = 247
$ws1.Range("F11").Value = 6

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 115
$ws4.Range("F4").Value = 500
$ws4.Range("F5").Value = 5063
$ws4.Range("F11").Value = 247
$ws4.Range("F12").Value = 6
